$wb = $excel.ActiveWorkbook

# Updated "想去人数" (want-to-go count) values for matching events across sheets.
# Sheet "展览" (sheetId 1) - rows keyed by row number
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 622
$ws1.Range("F4").Value = 607
$ws1.Range("F5").Value = 543
$ws1.Range("F6").Value = 303
$ws1.Range("F7").Value = 2727
$ws1.Range("F9").Value = 7612
$ws1.Range("F10").Value = 198
$ws1.Range("F11").Value = 464
$ws1.Range("F12").Value = 33
$ws1.Range("F13").Value = 282

# Sheet "全部类型" (sheetId 4) - same events, different row offsets due to extra
# "演出" rows mixed in
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 622
$ws4.Range("F4").Value = 607
$ws4.Range("F5").Value = 543
$ws4.Range("F6").Value = 303
$ws4.Range("F9").Value = 2727
$ws4.Range("F11").Value = 7612
$ws4.Range("F12").Value = 198
$ws4.Range("F13").Value = 464
$ws4.Range("F14").Value = 33
$ws4.Range("F17").Value = 282
